$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(31, 1).Value = 10001
$ws.Cells.Item(31, 2).Value = 10030
$ws.Cells.Item(31, 3).Value = "eng"
$ws.Cells.Item(31, 4).Value = $true
$ws.Cells.Item(31, 5).Value = "superadmin"
$ws.Cells.Item(31, 6).Value = "now()"

$ws.Cells.Item(32, 1).Value = 10001
$ws.Cells.Item(32, 2).Value = 10031
$ws.Cells.Item(32, 3).Value = "eng"
$ws.Cells.Item(32, 4).Value = $true
$ws.Cells.Item(32, 5).Value = "superadmin"
$ws.Cells.Item(32, 6).Value = "now()"

$ws.Range("E31").Select()
